$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.386.07"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "1.875.33"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7124"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.22"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3117"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07800"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08455"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").Value = "1.878.02"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.245"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7129"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.26"

$ws.Range("D16").Value = "29.387.73"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008253"
$ws.Range("E17").Value = "  +5.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.051"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.15"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("E20").Value = "  +0.93%  "

$ws.Range("D21").Value = "2.120.82"
$ws.Range("E21").Value = "  -1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.789"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1604"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.02"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.071"
$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.432"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("E31").Value = "  -4.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.318"
$ws.Range("E32").Value = "  +1.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05303"
$ws.Range("E33").Value = "  +2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.943"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.180"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7475"
$ws.Range("E36").Value = "  -6.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01874"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").Value = "1.220.99"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.477"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.98"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8878"
$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.78"
$ws.Range("E44").Value = "  +7.54%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "2.017.47"
$ws.Range("E46").Value = "  -1.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.822"
$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5209"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4324"
$ws.Range("E51").Value = "  +0.97%  "
